# Atualização de bases das ligas, do dia: 13-06-2024 às 19:35
# Rotate/swap the betting-odds data among rows 134/136/137 and 144/145
# (match ids in column A and the Div/Date in C/D stay put; every other
# column's value moves to a different row in the cycle described by the
# diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($Row, $B, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T, $U, $V, $W, $X, $Y, $Z, $AA, $AB, $AC, $AD) {
    $ws.Range("B$Row").Value = $B
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Value = $G
    $ws.Range("H$Row").Value = $H
    $ws.Range("I$Row").Value = $I
    $ws.Range("J$Row").Value = $J
    $ws.Range("K$Row").Value = $K
    $ws.Range("L$Row").Value = $L
    $ws.Range("M$Row").Value = $M
    $ws.Range("N$Row").Value = $N
    $ws.Range("O$Row").Value = $O
    $ws.Range("P$Row").Value = $P
    $ws.Range("Q$Row").Value = $Q
    $ws.Range("R$Row").Value = $R
    $ws.Range("S$Row").Value = $S
    $ws.Range("T$Row").Value = $T
    $ws.Range("U$Row").Value = $U
    $ws.Range("V$Row").Value = $V
    $ws.Range("W$Row").Value = $W
    $ws.Range("X$Row").Value = $X
    $ws.Range("Y$Row").Value = $Y
    $ws.Range("Z$Row").Value = $Z
    $ws.Range("AA$Row").Value = $AA
    $ws.Range("AB$Row").Value = $AB
    $ws.Range("AC$Row").Value = $AC
    $ws.Range("AD$Row").Value = $AD
}

# Row 134 <- data that used to live on row 136
Set-Row 134 7483306 "Tecnico Universitario" "Club Atletico Libertad" 1 1 0 0 "D" 1.5 4.333 5.75 1.533 4.2 5.5 -1 1.925 1.875 2.25 1.8 2 -1 3.2 -1 -1 0.875 -0.5 0.5

# Row 136 <- data that used to live on row 137
Set-Row 136 7482867 "Cumbaya FC" "LDU Quito" 1 2 0 0 "A" 5.25 3.75 1.65 9 4.5 1.363 1.25 1.975 1.825 2.5 1.825 1.975 -1 -1 0.363 0.4875 -0.5 0.825 -1

# Row 137 <- data that used to live on row 134
Set-Row 137 7483188 "Gualaceo SC" "Emelec" 0 2 0 1 "A" 3.6 3.3 2.05 2.6 3.25 2.75 0 1.8 2 2.5 1.975 1.825 -1 -1 1.75 -1 1 -1 0.825

# Row 144 <- data that used to live on row 145
Set-Row 144 7528857 "Universidad Catolica del Ecuador" "Barcelona Guayaquil" 0 1 0 0 "A" 1.533 4 5.5 1.5 4.333 5.25 -1 1.8 2 3 1.975 1.825 -1 -1 4.25 -1 1 -1 0.825

# Row 145 <- data that used to live on row 144
Set-Row 145 7528852 "Delfin SC" "Tecnico Universitario" 2 2 1 0 "D" 2.1 3.4 3.1 2.1 3.4 3.1 -0.25 1.8 2 2.25 1.9 1.9 -1 2.4 -1 -0.5 0.5 0.8999999999999999 -1

Write-Output "rows updated"
